$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column D
$ws.Range("D1").Value = "target01"

# Add new values for column D, rows 2-8
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0

# Update selection to D9, matching the diff
$ws.Range("D9").Select()
